# Scheduled-runner refresh of computed market/profit figures (columns H-N)
# across several rows on each sheet. Values below are the new snapshot
# pulled by the runner; row/column layout is unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 228
$ws.Range("I12").Value = 228
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 228
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -58
$ws.Range("N12").Value = ""

$ws.Range("H129").Value = 835.35
$ws.Range("I129").Value = 210.88889
$ws.Range("J129").Value = 1016.64514
$ws.Range("K129").Value = 632.6666700000001
$ws.Range("L129").Value = 3049.93542
$ws.Range("M129").Value = 4367.333329999999
$ws.Range("N129").Value = -13049.93542

$ws.Range("H132").Value = 21543832
$ws.Range("I132").Value = 25775372
$ws.Range("J132").Value = 1454.1818
$ws.Range("K132").Value = 77326116
$ws.Range("L132").Value = 4362.5454
$ws.Range("M132").Value = -77323586
$ws.Range("N132").Value = -9422.545399999999

$ws.Range("H138").Value = 1340.94
$ws.Range("I138").Value = 721.1017000000001
$ws.Range("J138").Value = 2232.9023
$ws.Range("K138").Value = 2163.3051
$ws.Range("L138").Value = 6698.706900000001
$ws.Range("M138").Value = 2976.6949
$ws.Range("N138").Value = -16978.7069

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1402.8

$ws.Range("H32").Value = 4253.21
$ws.Range("I32").Value = 3601.045
$ws.Range("J32").Value = 9529.817999999999
$ws.Range("K32").Value = 3601.045
$ws.Range("L32").Value = 9529.817999999999
$ws.Range("M32").Value = -3314.045
$ws.Range("N32").Value = -10103.818

$ws.Range("H74").Value = 808.325
$ws.Range("I74").Value = 610.2381
$ws.Range("J74").Value = 1027.2632
$ws.Range("K74").Value = 610.2381
$ws.Range("L74").Value = 1027.2632
$ws.Range("M74").Value = 263.7619
$ws.Range("N74").Value = -2775.2632

$ws.Range("H77").Value = 808.325
$ws.Range("I77").Value = 610.2381
$ws.Range("J77").Value = 1027.2632
$ws.Range("K77").Value = 3051.1905
$ws.Range("L77").Value = 5136.316000000001
$ws.Range("M77").Value = 1316.8095
$ws.Range("N77").Value = -13872.316

$ws.Range("H110").Value = 594.1111
$ws.Range("I110").Value = 595.30304
$ws.Range("J110").Value = 590.8333
$ws.Range("K110").Value = 595.30304
$ws.Range("L110").Value = 590.8333
$ws.Range("M110").Value = 1449.69696
$ws.Range("N110").Value = -4680.8333

$ws.Range("H132").Value = 2661015.2
$ws.Range("I132").Value = 3677710.8
$ws.Range("J132").Value = 1965.4615
$ws.Range("K132").Value = 11033132.4
$ws.Range("L132").Value = 5896.3845
$ws.Range("M132").Value = -11030602.4
$ws.Range("N132").Value = -10956.3845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 629.8461
$ws.Range("I94").Value = 798.3333
$ws.Range("J94").Value = 485.42856
$ws.Range("K94").Value = 798.3333
$ws.Range("L94").Value = 485.42856
$ws.Range("M94").Value = -347.3333
$ws.Range("N94").Value = -1387.42856

$ws.Range("H134").Value = 7764829.5
$ws.Range("I134").Value = 8786228
$ws.Range("J134").Value = 2200
$ws.Range("K134").Value = 26358684
$ws.Range("L134").Value = 6600
$ws.Range("M134").Value = -26356149
$ws.Range("N134").Value = -11670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5182.7363
$ws.Range("I31").Value = 903.807
$ws.Range("J31").Value = 21442.666
$ws.Range("K31").Value = 903.807
$ws.Range("L31").Value = 21442.666
$ws.Range("M31").Value = -608.807
$ws.Range("N31").Value = -22032.666

$ws.Range("H34").Value = 5182.7363
$ws.Range("I34").Value = 903.807
$ws.Range("J34").Value = 21442.666
$ws.Range("K34").Value = 903.807
$ws.Range("L34").Value = 21442.666
$ws.Range("M34").Value = -701.807
$ws.Range("N34").Value = -21846.666

$ws.Range("H58").Value = 3000167
$ws.Range("I58").Value = 3888910
$ws.Range("J58").Value = 10759.454
$ws.Range("K58").Value = 3888910
$ws.Range("L58").Value = 10759.454
$ws.Range("M58").Value = -3888707
$ws.Range("N58").Value = -11165.454

$ws.Range("H132").Value = 4905072.5
$ws.Range("I132").Value = 7093642
$ws.Range("J132").Value = 6844.7617
$ws.Range("K132").Value = 21280926
$ws.Range("L132").Value = 20534.2851
$ws.Range("M132").Value = -21278396
$ws.Range("N132").Value = -25594.2851

$ws.Range("H134").Value = 27345004
$ws.Range("I134").Value = 35715336
$ws.Range("J134").Value = 4809502
$ws.Range("K134").Value = 107146008
$ws.Range("L134").Value = 14428506
$ws.Range("M134").Value = -107143473
$ws.Range("N134").Value = -14433576

$ws.Range("H136").Value = 3000167
$ws.Range("I136").Value = 3888910
$ws.Range("J136").Value = 10759.454
$ws.Range("K136").Value = 11666730
$ws.Range("L136").Value = 32278.362
$ws.Range("M136").Value = -11664180
$ws.Range("N136").Value = -37378.362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 430.5484
$ws.Range("I5").Value = 297.35
$ws.Range("J5").Value = 672.7273
$ws.Range("K5").Value = 892.0500000000001
$ws.Range("L5").Value = 2018.1819
$ws.Range("M5").Value = -780.0500000000001
$ws.Range("N5").Value = -2242.1819

$ws.Range("H9").Value = 500450000
$ws.Range("J9").Value = 900000
$ws.Range("L9").Value = 2700000
$ws.Range("N9").Value = -2700448

$ws.Range("H11").Value = 794.5
$ws.Range("I11").Value = 103.833336
$ws.Range("J11").Value = 2175.8333
$ws.Range("K11").Value = 311.500008
$ws.Range("L11").Value = 6527.499899999999
$ws.Range("M11").Value = -171.500008
$ws.Range("N11").Value = -6807.499899999999

$ws.Range("H113").Value = 5561666.5
$ws.Range("J113").Value = 12513001
$ws.Range("L113").Value = 37539003
$ws.Range("N113").Value = -37543343

$ws.Range("H122").Value = 794.2121
$ws.Range("I122").Value = 773.5
$ws.Range("J122").Value = 871.1429000000001
$ws.Range("K122").Value = 6961.5
$ws.Range("L122").Value = 7840.2861
$ws.Range("M122").Value = -4511.5
$ws.Range("N122").Value = -12740.2861

$ws.Range("H135").Value = 430.5484
$ws.Range("I135").Value = 297.35
$ws.Range("J135").Value = 672.7273
$ws.Range("K135").Value = 2676.15
$ws.Range("L135").Value = 6054.545700000001
$ws.Range("M135").Value = -141.1500000000001
$ws.Range("N135").Value = -11124.5457

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3847793.5
$ws.Range("I132").Value = 5000971.5
$ws.Range("J132").Value = 3867.3333
$ws.Range("K132").Value = 15002914.5
$ws.Range("L132").Value = 11601.9999
$ws.Range("M132").Value = -15000384.5
$ws.Range("N132").Value = -16661.9999

$ws.Range("H136").Value = 6367.607
$ws.Range("I136").Value = 8571.556
$ws.Range("J136").Value = 2400.5
$ws.Range("K136").Value = 25714.668
$ws.Range("L136").Value = 7201.5
$ws.Range("M136").Value = -23164.668
$ws.Range("N136").Value = -12301.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 125000660
$ws.Range("I4").Value = 699.8
$ws.Range("J4").Value = 333333920
$ws.Range("K4").Value = 699.8
$ws.Range("L4").Value = 333333920
$ws.Range("M4").Value = -586.8
$ws.Range("N4").Value = -333334146

$ws.Range("H92").Value = 39000
$ws.Range("J92").Value = 39000
$ws.Range("L92").Value = 39000
$ws.Range("N92").Value = -43992
